# Remove the "Mô tả" (Description) column (column C) from the product sample sheet.
# This shifts columns D:G left to C:F, matching the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C:C").Delete()

$ws.Range("G10").Select()

$wb.Save()
